# Updated cryptos list on Fri Jun 21 08:18:44 UTC 2024 with GitHub Actions
#
# Applies the latest coinranking.com price/volume snapshot to Sheet1,
# and reorders two pairs of coins (ImmutableX/Aptos and
# FirstDigitalUSD/EnergySwap) whose ranking changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing the cell to remain plain text
# (keeps e.g. "73.00" / "0.999" / "581.99" from being re-interpreted as
# numbers, which would silently drop trailing zeros). The NumberFormat
# is reset back to the default "Normal" style right after, so no stray
# formatting is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "64.300.85"
$ws.Range("E2").Value = "  -2.22%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "3.482.09"
$ws.Range("E3").Value = "  -3.32%  "

# --- Row 4: TetherUSD ---
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.11%  "

# --- Row 5: BNB ---
Set-TextValue $ws.Range("D5") "581.99"
$ws.Range("E5").Value = "  -3.61%  "

# --- Row 6: Solana ---
Set-TextValue $ws.Range("D6") "131.32"
$ws.Range("E6").Value = "  -4.07%  "

# --- Row 7: LidoStakedEther ---
$ws.Range("D7").Value = "3.481.85"
$ws.Range("E7").Value = "  -3.30%  "

# --- Row 8: USDC ---
$ws.Range("E8").Value = "  +0.01%  "

# --- Row 9: XRP ---
Set-TextValue $ws.Range("D9") "0.490"
$ws.Range("E9").Value = "  -1.88%  "

# --- Row 10: Dogecoin ---
$ws.Range("E10").Value = "  -1.49%  "

# --- Row 11: Toncoin ---
Set-TextValue $ws.Range("D11") "7.22"
$ws.Range("E11").Value = "  -0.18%  "

# --- Row 12: Cardano ---
Set-TextValue $ws.Range("D12") "0.385"
$ws.Range("E12").Value = "  -1.66%  "

# --- Row 13: WrappedliquidstakedEther2.0 ---
$ws.Range("D13").Value = "4.069.60"
$ws.Range("E13").Value = "  -3.55%  "

# --- Row 14: Avalanche ---
Set-TextValue $ws.Range("D14") "27.62"
$ws.Range("E14").Value = "  -1.93%  "

# --- Row 15: ShibaInu ---
$ws.Range("E15").Value = "  -4.88%  "

# --- Row 16: TRON ---
$ws.Range("E16").Value = "  +0.34%  "

# --- Row 17: WrappedEther ---
$ws.Range("D17").Value = "3.477.41"
$ws.Range("E17").Value = "  -3.52%  "

# --- Row 18: WrappedBTC ---
$ws.Range("D18").Value = "64.307.55"
$ws.Range("E18").Value = "  -2.40%  "

# --- Row 19: Uniswap ---
Set-TextValue $ws.Range("D19") "9.87"
$ws.Range("E19").Value = "  -2.04%  "

# --- Row 20: Chainlink ---
Set-TextValue $ws.Range("D20") "14.27"
$ws.Range("E20").Value = "  -3.04%  "

# --- Row 21: Polkadot ---
$ws.Range("E21").Value = "  -4.38%  "

# --- Row 22: BitcoinCash ---
Set-TextValue $ws.Range("D22") "391.06"
$ws.Range("E22").Value = "  -1.67%  "

# --- Row 23: Polygon ---
$ws.Range("E23").Value = "  -2.93%  "

# --- Row 24: WrappedeETH ---
$ws.Range("D24").Value = "3.621.14"
$ws.Range("E24").Value = "  -3.41%  "

# --- Row 25: Litecoin ---
Set-TextValue $ws.Range("D25") "73.00"
$ws.Range("E25").Value = "  -2.09%  "

# --- Row 26: Dai ---
$ws.Range("E26").Value = "  +0.32%  "

# --- Row 27: PEPE ---
$ws.Range("E27").Value = "  -8.80%  "

# --- Row 28: Fetch.AI ---
Set-TextValue $ws.Range("D28") "1.55"
$ws.Range("E28").Value = "  -6.71%  "

# --- Row 29: Binance-PegBSC-USD ---
$ws.Range("E29").Value = "  +0.18%  "

# --- Row 30: RenderToken ---
Set-TextValue $ws.Range("D30") "7.36"
$ws.Range("E30").Value = "  -9.78%  "

# --- Row 31: PancakeSwap ---
Set-TextValue $ws.Range("D31") "2.25"
$ws.Range("E31").Value = "  -6.94%  "

# --- Row 32: InternetComputer(DFINITY) ---
Set-TextValue $ws.Range("D32") "8.16"
$ws.Range("E32").Value = "  -5.21%  "

# --- Row 33: RenzoRestakedETH ---
$ws.Range("D33").Value = "3.477.77"
$ws.Range("E33").Value = "  -3.51%  "

# --- Row 34: USDe ---
$ws.Range("E34").Value = "  +0.02%  "

# --- Row 35: EthereumClassic ---
Set-TextValue $ws.Range("D35") "23.82"
$ws.Range("E35").Value = "  -3.01%  "

# --- Row 36: Kaspa ---
Set-TextValue $ws.Range("D36") "0.144"
$ws.Range("E36").Value = "  -2.70%  "

# --- Row 37: NEARProtocol ---
Set-TextValue $ws.Range("D37") "5.20"
$ws.Range("E37").Value = "  -3.65%  "

# --- Row 38: Monero ---
Set-TextValue $ws.Range("D38") "169.97"
$ws.Range("E38").Value = "  -0.38%  "

# --- Rows 39 & 40: ImmutableX and Aptos swap ranking positions ---
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D39") "1.57"
$ws.Range("E39").Value = "  -2.70%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D40") "6.95"
$ws.Range("E40").Value = "  -2.00%  "

# --- Row 41: Hedera ---
Set-TextValue $ws.Range("D41") "0.0802"
$ws.Range("E41").Value = "  -4.24%  "

# --- Row 42: Mantle ---
Set-TextValue $ws.Range("D42") "0.809"
$ws.Range("E42").Value = "  -3.93%  "

# --- Rows 43 & 44: FirstDigitalUSD and EnergySwap swap ranking positions ---
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D43") "0.999"
$ws.Range("E43").Value = "  -0.19%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "25.41"
$ws.Range("E44").Value = "  -3.30%  "

# --- Row 45: OKB ---
Set-TextValue $ws.Range("D45") "41.69"
$ws.Range("E45").Value = "  -3.79%  "

# --- Row 46: ONDO ---
Set-TextValue $ws.Range("D46") "1.19"
$ws.Range("E46").Value = "  -5.21%  "

# --- Row 47: Filecoin ---
$ws.Range("E47").Value = "  -4.23%  "

# --- Row 48: Stacks ---
$ws.Range("E48").Value = "  -4.39%  "

# --- Row 49: Cosmos ---
$ws.Range("E49").Value = "  -3.12%  "

# --- Row 50: Maker ---
$ws.Range("D50").Value = "2.428.85"
$ws.Range("E50").Value = "  +0.01%  "

# --- Row 51: SuiNetwork ---
Set-TextValue $ws.Range("D51") "0.887"
$ws.Range("E51").Value = "  -0.99%  "
